# Commit: "chinh sua lan 1 / ok hoan tat" -- append " da chinh sua" (Vietnamese
# for "already edited") as a new run right after the existing "*font-size:"
# run in the first paragraph, and move the hidden "_GoBack" bookmark (which
# Word keeps pinned to the location of the most recent edit) from the end of
# the "- bottom:" paragraph to the end of this now-edited paragraph.

$d = $word.ActiveDocument

# --- 1. Locate the first paragraph ("*font-size:") and append new text ----
$p1 = $d.Paragraphs(1)

$r = $p1.Range
[void]$r.MoveEnd(1, -1)    # exclude the paragraph mark
$insertionPoint = $r.End   # offset right after "*font-size:"

# Collapse to the end and insert the new text. A one-character placeholder
# ("Z") is appended too -- it gives us a safe, non-boundary anchor point to
# park the _GoBack bookmark on in step 3, and is deleted again afterwards.
$tail = $r.Duplicate
$tail.Collapse(0)
$tail.InsertAfter(" đã chĩnh sửa" + "Z")

# --- 2. Force the new text into its own run (matching formatting) ---------
# InsertAfter merges into the preceding run if formatting is identical, but
# the target document has the new text as a *separate* <w:r>. Toggling the
# (identical) font size on just the inserted span causes the engine to split
# the run without changing the resulting formatting.
$p1b = $d.Paragraphs(1)
$paraEnd = $p1b.Range.End          # end, including the paragraph mark
$newSpanEnd = $paraEnd - 1         # end, excluding the paragraph mark
$newSpan = $d.Range($insertionPoint, $newSpanEnd)
$newSpan.Font.Size = 20
$newSpan.Font.Size = 18

# --- 3. Move the _GoBack bookmark to the end of this paragraph ------------
# Re-anchor it on the trailing "Z" placeholder (a safe, non-collapsed range),
# then shrink it back to zero width by clearing that one character through
# the bookmark's own range -- this leaves a correctly collapsed bookmark
# right where Word would leave _GoBack after the most recent edit.
$placeholder = $d.Range($newSpanEnd - 1, $newSpanEnd)

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$d.Bookmarks.Add("_GoBack", $placeholder)

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Range.Text = ""
